# Insert a new weekly pair of rows (Primera/Segunda) before the existing
# row 44, shifting all subsequent data rows down by two. Then populate the
# two newly-inserted rows with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 44.. down by 2 (inserts two blank rows at 44 and 45, copying
# formatting from the row above as Excel normally does on row insert).
$ws.Rows("44:45").Insert()

# Populate the newly inserted row 44 ("Primera" quality) with this week's data.
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C44").Value = "Arica y Parinacota"
$ws.Range("D44").Value = 44481
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = 100114014
$ws.Range("G44").Value = "Betarraga"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 1200
$ws.Range("K44").Value = 400
$ws.Range("L44").Value = 450
$ws.Range("M44").Value = 425
$ws.Range("N44").Value = "`$/paquete 4 unidades"
$ws.Range("O44").Value = "Región de Arica y Parinacota"
$ws.Range("P44").Value = 106
$ws.Range("Q44").Value = 4
$ws.Range("R44").Value = "Hortaliza"

# Populate the newly inserted row 45 ("Segunda" quality) with this week's data.
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C45").Value = "Arica y Parinacota"
$ws.Range("D45").Value = 44481
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 100114014
$ws.Range("G45").Value = "Betarraga"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Segunda"
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 400
$ws.Range("L45").Value = 450
$ws.Range("M45").Value = 425
$ws.Range("N45").Value = "`$/paquete 5 unidades"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value = 85
$ws.Range("Q45").Value = 5
$ws.Range("R45").Value = "Hortaliza"
